$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the rest of row 1 (e.g. H1) by copying
# its formatting onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for the new I (I0) and J (IF) columns, rows 2..29
$data = @{
    2  = @(1, 5)
    3  = @(1, 6)
    4  = @(1, 4)
    5  = @(1, 5)
    6  = @(1, 6)
    7  = @(1, 6)
    8  = @(1, 7)
    9  = @(1, 5)
    10 = @(1, 5)
    11 = @(1, 7)
    12 = @(1, 5)
    13 = @(1, 5)
    14 = @(1, 3)
    15 = @(4, 6)
    16 = @(9, 9)
    17 = @(8, 8)
    18 = @(6, 7)
    19 = @(3, 5)
    20 = @(3, 5)
    21 = @(3, 6)
    22 = @(6, 7)
    23 = @(7, 8)
    24 = @(8, 8)
    25 = @(8, 8)
    26 = @(4, 4)
    27 = @(2, 3)
    28 = @(2, 3)
    29 = @(2, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
